# Auto-generated edits applying numeric updates to Leve profit sheets (Behemoth_Profits)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3831.818   # H40: 3784.5386 -> 3831.818
$ws.Cells.Item(40, 9).Value = 3525   # I40: 3524.875 -> 3525
$ws.Cells.Item(40, 11).Value = 3525   # K40: 3524.875 -> 3525
$ws.Cells.Item(40, 13).Value = -3350   # M40: -3349.875 -> -3350
$ws.Cells.Item(62, 8).Value = 16749.75   # H62: 20500 -> 16749.75
$ws.Cells.Item(62, 9).Value = 4499.5   # I62: 6000 -> 4499.5
$ws.Cells.Item(62, 10).Value = 20833.166   # J62: 23400 -> 20833.166
$ws.Cells.Item(62, 11).Value = 4499.5   # K62: 6000 -> 4499.5
$ws.Cells.Item(62, 12).Value = 20833.166   # L62: 23400 -> 20833.166
$ws.Cells.Item(62, 13).Value = -3875.5   # M62: -5376 -> -3875.5
$ws.Cells.Item(62, 14).Value = -22081.166   # N62: -24648 -> -22081.166
$ws.Cells.Item(65, 8).Value = 16749.75   # H65: 20500 -> 16749.75
$ws.Cells.Item(65, 9).Value = 4499.5   # I65: 6000 -> 4499.5
$ws.Cells.Item(65, 10).Value = 20833.166   # J65: 23400 -> 20833.166
$ws.Cells.Item(65, 11).Value = 22497.5   # K65: 30000 -> 22497.5
$ws.Cells.Item(65, 12).Value = 104165.83   # L65: 117000 -> 104165.83
$ws.Cells.Item(65, 13).Value = -19377.5   # M65: -26880 -> -19377.5
$ws.Cells.Item(65, 14).Value = -110405.83   # N65: -123240 -> -110405.83
$ws.Cells.Item(111, 8).Value = 4171.3335   # H111: 3532 -> 4171.3335
$ws.Cells.Item(111, 9).Value = 4205.6   # I111: 3465.1428 -> 4205.6
$ws.Cells.Item(111, 11).Value = 12616.8   # K111: 10395.4284 -> 12616.8
$ws.Cells.Item(111, 13).Value = -9549.800000000001   # M111: -7328.428400000001 -> -9549.800000000001
$ws.Cells.Item(137, 8).Value = 7450.864   # H137: 7296.478 -> 7450.864
$ws.Cells.Item(137, 9).Value = 6943.263   # I137: 6791.1 -> 6943.263
$ws.Cells.Item(137, 11).Value = 20829.789   # K137: 20373.3 -> 20829.789
$ws.Cells.Item(137, 13).Value = -18279.789   # M137: -17823.3 -> -18279.789
$ws.Cells.Item(141, 8).Value = 10527.357   # H141: 11299.417 -> 10527.357
$ws.Cells.Item(141, 9).Value = 8917.5   # I141: 9673.125 -> 8917.5
$ws.Cells.Item(141, 11).Value = 26752.5   # K141: 29019.375 -> 26752.5
$ws.Cells.Item(141, 13).Value = -21572.5   # M141: -23839.375 -> -21572.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 455.3   # H5: 489.8889 -> 455.3
$ws.Cells.Item(5, 9).Value = 162.5   # I5: 168.66667 -> 162.5
$ws.Cells.Item(5, 11).Value = 162.5   # K5: 168.66667 -> 162.5
$ws.Cells.Item(5, 13).Value = -50.5   # M5: -56.66667000000001 -> -50.5
$ws.Cells.Item(7, 8).Value = 90542.39999999999   # H7: 93750 -> 90542.39999999999
$ws.Cells.Item(7, 10).Value = 97570.664   # J7: 107500 -> 97570.664
$ws.Cells.Item(7, 12).Value = 97570.664   # L7: 107500 -> 97570.664
$ws.Cells.Item(7, 14).Value = -97798.664   # N7: -107728 -> -97798.664
$ws.Cells.Item(32, 8).Value = 15156092   # H32: 15154900 -> 15156092
$ws.Cells.Item(32, 9).Value = 16669034   # I32: 15627240 -> 16669034
$ws.Cells.Item(32, 10).Value = 26671.334   # J32: 40014 -> 26671.334
$ws.Cells.Item(32, 11).Value = 16669034   # K32: 15627240 -> 16669034
$ws.Cells.Item(32, 12).Value = 26671.334   # L32: 40014 -> 26671.334
$ws.Cells.Item(32, 13).Value = -16668747   # M32: -15626953 -> -16668747
$ws.Cells.Item(32, 14).Value = -27245.334   # N32: -40588 -> -27245.334
$ws.Cells.Item(132, 8).Value = 3188.4   # H132: 2788.2974 -> 3188.4
$ws.Cells.Item(132, 9).Value = 3160.4482   # I132: 2964.1562 -> 3160.4482
$ws.Cells.Item(132, 10).Value = 3999   # J132: 1662.8 -> 3999
$ws.Cells.Item(132, 11).Value = 9481.3446   # K132: 8892.4686 -> 9481.3446
$ws.Cells.Item(132, 12).Value = 11997   # L132: 4988.4 -> 11997
$ws.Cells.Item(132, 13).Value = -6951.3446   # M132: -6362.4686 -> -6951.3446
$ws.Cells.Item(132, 14).Value = -17057   # N132: -10048.4 -> -17057

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 455.3   # H4: 489.8889 -> 455.3
$ws.Cells.Item(4, 9).Value = 162.5   # I4: 168.66667 -> 162.5
$ws.Cells.Item(4, 11).Value = 162.5   # K4: 168.66667 -> 162.5
$ws.Cells.Item(4, 13).Value = -47.5   # M4: -53.66667000000001 -> -47.5
$ws.Cells.Item(81, 8).Value = 28405.25   # H81: 30155 -> 28405.25
$ws.Cells.Item(81, 10).Value = 28405.25   # J81: 30155 -> 28405.25
$ws.Cells.Item(81, 12).Value = 28405.25   # L81: 30155 -> 28405.25
$ws.Cells.Item(81, 14).Value = -30527.25   # N81: -32277 -> -30527.25
$ws.Cells.Item(84, 8).Value = 28405.25   # H84: 30155 -> 28405.25
$ws.Cells.Item(84, 10).Value = 28405.25   # J84: 30155 -> 28405.25
$ws.Cells.Item(84, 12).Value = 85215.75   # L84: 90465 -> 85215.75
$ws.Cells.Item(84, 14).Value = -95823.75   # N84: -101073 -> -95823.75
$ws.Cells.Item(105, 8).Value = 2462.4348   # H105: 2399.4167 -> 2462.4348
$ws.Cells.Item(105, 9).Value = 1134.1666   # I105: 1107.8572 -> 1134.1666
$ws.Cells.Item(105, 11).Value = 1134.1666   # K105: 1107.8572 -> 1134.1666
$ws.Cells.Item(105, 13).Value = 612.8334   # M105: 639.1428000000001 -> 612.8334
$ws.Cells.Item(134, 8).Value = 582355.5   # H134: 480351.75 -> 582355.5
$ws.Cells.Item(134, 9).Value = 1885.6364   # I134: 1530.3684 -> 1885.6364
$ws.Cells.Item(134, 10).Value = 1380501.5   # J134: 2754753.2 -> 1380501.5
$ws.Cells.Item(134, 11).Value = 5656.9092   # K134: 4591.1052 -> 5656.9092
$ws.Cells.Item(134, 12).Value = 4141504.5   # L134: 8264259.600000001 -> 4141504.5
$ws.Cells.Item(134, 13).Value = -3121.9092   # M134: -2056.1052 -> -3121.9092
$ws.Cells.Item(134, 14).Value = -4146574.5   # N134: -8269329.600000001 -> -4146574.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(45, 8).Value = 0   # H45: 8995 -> 0
$ws.Cells.Item(45, 10).Value = 0   # J45: 8995 -> 0
$ws.Cells.Item(45, 12).Value = 0   # L45: 8995 -> 0
$ws.Cells.Item(45, 14).Value = ""   # N45: clear (was -10181)
$ws.Cells.Item(99, 8).Value = 3100.375   # H99: 3022.6667 -> 3100.375
$ws.Cells.Item(99, 9).Value = 2979.6   # I99: 2979.8 -> 2979.6
$ws.Cells.Item(99, 10).Value = 3301.6667   # J99: 3076.25 -> 3301.6667
$ws.Cells.Item(99, 11).Value = 2979.6   # K99: 2979.8 -> 2979.6
$ws.Cells.Item(99, 12).Value = 3301.6667   # L99: 3076.25 -> 3301.6667
$ws.Cells.Item(99, 13).Value = -1481.6   # M99: -1481.8 -> -1481.6
$ws.Cells.Item(99, 14).Value = -6297.6667   # N99: -6072.25 -> -6297.6667
$ws.Cells.Item(110, 8).Value = 103344.5   # H110: 87795 -> 103344.5
$ws.Cells.Item(110, 10).Value = 103344.5   # J110: 87795 -> 103344.5
$ws.Cells.Item(110, 12).Value = 103344.5   # L110: 87795 -> 103344.5
$ws.Cells.Item(110, 14).Value = -111524.5   # N110: -95975 -> -111524.5
$ws.Cells.Item(111, 8).Value = 99989   # H111: 99845.5 -> 99989
$ws.Cells.Item(111, 10).Value = 99989   # J111: 99845.5 -> 99989
$ws.Cells.Item(111, 12).Value = 99989   # L111: 99845.5 -> 99989
$ws.Cells.Item(111, 14).Value = -108169   # N111: -108025.5 -> -108169
$ws.Cells.Item(124, 8).Value = 57435   # H124: 42160.75 -> 57435
$ws.Cells.Item(124, 10).Value = 57435   # J124: 42160.75 -> 57435
$ws.Cells.Item(124, 12).Value = 57435   # L124: 42160.75 -> 57435
$ws.Cells.Item(124, 14).Value = -62345   # N124: -47070.75 -> -62345
$ws.Cells.Item(126, 8).Value = 3100.375   # H126: 3022.6667 -> 3100.375
$ws.Cells.Item(126, 9).Value = 2979.6   # I126: 2979.8 -> 2979.6
$ws.Cells.Item(126, 10).Value = 3301.6667   # J126: 3076.25 -> 3301.6667
$ws.Cells.Item(126, 11).Value = 8938.799999999999   # K126: 8939.400000000001 -> 8938.799999999999
$ws.Cells.Item(126, 12).Value = 9905.000100000001   # L126: 9228.75 -> 9905.000100000001
$ws.Cells.Item(126, 13).Value = -6468.799999999999   # M126: -6469.400000000001 -> -6468.799999999999
$ws.Cells.Item(126, 14).Value = -14845.0001   # N126: -14168.75 -> -14845.0001
$ws.Cells.Item(134, 8).Value = 3140.2942   # H134: 2403.2 -> 3140.2942
$ws.Cells.Item(134, 9).Value = 2706.0715   # I134: 2026.3636 -> 2706.0715
$ws.Cells.Item(134, 11).Value = 8118.2145   # K134: 6079.0908 -> 8118.2145
$ws.Cells.Item(134, 13).Value = -5583.2145   # M134: -3544.0908 -> -5583.2145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1438.4546   # H113: 1402.5555 -> 1438.4546
$ws.Cells.Item(113, 10).Value = 1502.4   # J113: 1478 -> 1502.4
$ws.Cells.Item(113, 12).Value = 4507.200000000001   # L113: 4434 -> 4507.200000000001
$ws.Cells.Item(113, 14).Value = -8847.200000000001   # N113: -8774 -> -8847.200000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3491.8948   # H102: 1963.8049 -> 3491.8948
$ws.Cells.Item(102, 9).Value = 2780   # I102: 1436.7742 -> 2780
$ws.Cells.Item(102, 10).Value = 4712.2856   # J102: 3597.6 -> 4712.2856
$ws.Cells.Item(102, 11).Value = 2780   # K102: 1436.7742 -> 2780
$ws.Cells.Item(102, 12).Value = 4712.2856   # L102: 3597.6 -> 4712.2856
$ws.Cells.Item(102, 13).Value = -1158   # M102: 185.2257999999999 -> -1158
$ws.Cells.Item(102, 14).Value = -7956.2856   # N102: -6841.6 -> -7956.2856
$ws.Cells.Item(113, 8).Value = 4427.7856   # H113: 3366.7144 -> 4427.7856
$ws.Cells.Item(113, 9).Value = 3498.5   # I113: 2021.3846 -> 3498.5
$ws.Cells.Item(113, 10).Value = 4582.6665   # J113: 4532.6665 -> 4582.6665
$ws.Cells.Item(113, 11).Value = 3498.5   # K113: 2021.3846 -> 3498.5
$ws.Cells.Item(113, 12).Value = 4582.6665   # L113: 4532.6665 -> 4582.6665
$ws.Cells.Item(113, 13).Value = -1328.5   # M113: 148.6153999999999 -> -1328.5
$ws.Cells.Item(113, 14).Value = -8922.666499999999   # N113: -8872.666499999999 -> -8922.666499999999
$ws.Cells.Item(126, 8).Value = 4503.2856   # H126: 4565.077 -> 4503.2856
$ws.Cells.Item(126, 9).Value = 4006.5715   # I126: 4057.6667 -> 4006.5715
$ws.Cells.Item(126, 11).Value = 12019.7145   # K126: 12173.0001 -> 12019.7145
$ws.Cells.Item(126, 13).Value = -9549.7145   # M126: -9703.000100000001 -> -9549.7145

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1137.5   # H16: 904.2857 -> 1137.5
$ws.Cells.Item(16, 9).Value = 1137.5   # I16: 904.2857 -> 1137.5
$ws.Cells.Item(16, 11).Value = 1137.5   # K16: 904.2857 -> 1137.5
$ws.Cells.Item(16, 13).Value = -967.5   # M16: -734.2857 -> -967.5
$ws.Cells.Item(22, 8).Value = 1627.4546   # H22: 2491.8333 -> 1627.4546
$ws.Cells.Item(22, 9).Value = 1700   # I22: 2433.3333 -> 1700
$ws.Cells.Item(22, 10).Value = 1500.5   # J22: 2667.3333 -> 1500.5
$ws.Cells.Item(22, 11).Value = 1700   # K22: 2433.3333 -> 1700
$ws.Cells.Item(22, 12).Value = 1500.5   # L22: 2667.3333 -> 1500.5
$ws.Cells.Item(22, 13).Value = -1405   # M22: -2138.3333 -> -1405
$ws.Cells.Item(22, 14).Value = -2090.5   # N22: -3257.3333 -> -2090.5
$ws.Cells.Item(24, 8).Value = 83340.336   # H24: 85007 -> 83340.336
$ws.Cells.Item(24, 10).Value = 83340.336   # J24: 85007 -> 83340.336
$ws.Cells.Item(24, 12).Value = 83340.336   # L24: 85007 -> 83340.336
$ws.Cells.Item(24, 14).Value = -84026.336   # N24: -85693 -> -84026.336
$ws.Cells.Item(27, 8).Value = 1627.4546   # H27: 2491.8333 -> 1627.4546
$ws.Cells.Item(27, 9).Value = 1700   # I27: 2433.3333 -> 1700
$ws.Cells.Item(27, 10).Value = 1500.5   # J27: 2667.3333 -> 1500.5
$ws.Cells.Item(27, 11).Value = 1700   # K27: 2433.3333 -> 1700
$ws.Cells.Item(27, 12).Value = 1500.5   # L27: 2667.3333 -> 1500.5
$ws.Cells.Item(27, 13).Value = -1593   # M27: -2326.3333 -> -1593
$ws.Cells.Item(27, 14).Value = -1714.5   # N27: -2881.3333 -> -1714.5
$ws.Cells.Item(46, 8).Value = 2176.4666   # H46: 1790.5 -> 2176.4666
$ws.Cells.Item(46, 9).Value = 1413.3636   # I46: 962.3333 -> 1413.3636
$ws.Cells.Item(46, 11).Value = 1413.3636   # K46: 962.3333 -> 1413.3636
$ws.Cells.Item(46, 13).Value = -1225.3636   # M46: -774.3333 -> -1225.3636
$ws.Cells.Item(61, 8).Value = 2820.3635   # H61: 2639.4167 -> 2820.3635
$ws.Cells.Item(61, 9).Value = 2317.1428   # I61: 2108.625 -> 2317.1428
$ws.Cells.Item(61, 11).Value = 2317.1428   # K61: 2108.625 -> 2317.1428
$ws.Cells.Item(61, 13).Value = -2115.1428   # M61: -1906.625 -> -2115.1428
$ws.Cells.Item(113, 8).Value = 2820.3635   # H113: 2639.4167 -> 2820.3635
$ws.Cells.Item(113, 9).Value = 2317.1428   # I113: 2108.625 -> 2317.1428
$ws.Cells.Item(113, 11).Value = 2317.1428   # K113: 2108.625 -> 2317.1428
$ws.Cells.Item(113, 13).Value = -147.1428000000001   # M113: 61.375 -> -147.1428000000001
$ws.Cells.Item(122, 8).Value = 6587.2085   # H122: 6702.227 -> 6587.2085
$ws.Cells.Item(122, 9).Value = 5926.5625   # I122: 6038.7334 -> 5926.5625
$ws.Cells.Item(122, 10).Value = 7908.5   # J122: 8124 -> 7908.5
$ws.Cells.Item(122, 11).Value = 17779.6875   # K122: 18116.2002 -> 17779.6875
$ws.Cells.Item(122, 12).Value = 23725.5   # L122: 24372 -> 23725.5
$ws.Cells.Item(122, 13).Value = -15329.6875   # M122: -15666.2002 -> -15329.6875
$ws.Cells.Item(122, 14).Value = -28625.5   # N122: -29272 -> -28625.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 0   # H3: 3999 -> 0
$ws.Cells.Item(3, 9).Value = 0   # I3: 3999 -> 0
$ws.Cells.Item(3, 11).Value = 0   # K3: 3999 -> 0
$ws.Cells.Item(3, 13).Value = ""   # M3: clear (was -3885)
$ws.Cells.Item(52, 8).Value = 22367   # H52: 19015.2 -> 22367
$ws.Cells.Item(52, 9).Value = 14984.25   # I52: 15019.25 -> 14984.25
$ws.Cells.Item(52, 10).Value = 29749.75   # J52: 34999 -> 29749.75
$ws.Cells.Item(52, 11).Value = 14984.25   # K52: 15019.25 -> 14984.25
$ws.Cells.Item(52, 12).Value = 29749.75   # L52: 34999 -> 29749.75
$ws.Cells.Item(52, 13).Value = -14758.25   # M52: -14793.25 -> -14758.25
$ws.Cells.Item(52, 14).Value = -30201.75   # N52: -35451 -> -30201.75
$ws.Cells.Item(124, 8).Value = 127994.5   # H124: 91996.336 -> 127994.5
$ws.Cells.Item(124, 10).Value = 127994.5   # J124: 91996.336 -> 127994.5
$ws.Cells.Item(124, 12).Value = 127994.5   # L124: 91996.336 -> 127994.5
$ws.Cells.Item(124, 14).Value = -137814.5   # N124: -101816.336 -> -137814.5
$ws.Cells.Item(132, 8).Value = 559183.8   # H132: 628869.7 -> 559183.8
$ws.Cells.Item(132, 9).Value = 5063.125   # I132: 5985.1665 -> 5063.125
$ws.Cells.Item(132, 10).Value = 1002480.4   # J132: 1002600.4 -> 1002480.4
$ws.Cells.Item(132, 11).Value = 15189.375   # K132: 17955.4995 -> 15189.375
$ws.Cells.Item(132, 12).Value = 3007441.2   # L132: 3007801.2 -> 3007441.2
$ws.Cells.Item(132, 13).Value = -12659.375   # M132: -15425.4995 -> -12659.375
$ws.Cells.Item(132, 14).Value = -3012501.2   # N132: -3012861.2 -> -3012501.2
